# Update countries & provincias Spain
# Applies the COVID-19 daily data refresh: country-name relabeling (shared-string
# content swaps) + updated case counters, and bumps the "datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name corrections (column A) ---
$ws.Range("A21").Value = 'Brasil'
$ws.Range("A22").Value = 'Israel'
$ws.Range("A144").Value = 'Congo'
$ws.Range("A145").Value = 'Etiopia'
$ws.Range("A147").Value = 'Mali'
$ws.Range("A148").Value = 'Niger'
$ws.Range("A156").Value = 'Eritrea'
$ws.Range("A157").Value = 'Guinea Ecuatorial'
$ws.Range("A159").Value = 'San Martin (Parte Francesa)'
$ws.Range("A160").Value = 'Bahamas'
$ws.Range("A161").Value = 'Namibia'
$ws.Range("A162").Value = 'Birmania'
$ws.Range("A163").Value = 'Groenlandia'
$ws.Range("A165").Value = 'Suazilandia'
$ws.Range("A166").Value = 'Siria'
$ws.Range("A167").Value = 'Laos'
$ws.Range("A168").Value = 'Seychelles'
$ws.Range("A169").Value = 'Surinam'
$ws.Range("A170").Value = 'Mozambique'
$ws.Range("A171").Value = 'Libia'
$ws.Range("A173").Value = 'Guyana'
$ws.Range("A174").Value = 'Curazao'
$ws.Range("A175").Value = 'Antigua y Barbuda'
$ws.Range("A176").Value = 'Zimbabue'
$ws.Range("A177").Value = 'Gabon'
$ws.Range("A178").Value = 'Angola'
$ws.Range("A179").Value = 'Santa Sede'
$ws.Range("A180").Value = 'San Martin (Parte Holandesa)'
$ws.Range("A182").Value = 'Cabo Verde'
$ws.Range("A183").Value = 'Sudan'
$ws.Range("A191").Value = 'Gambia'
$ws.Range("A192").Value = 'Nicaragua'
$ws.Range("A193").Value = 'Santa Lucia'
$ws.Range("A194").Value = 'Republica de Africa Central'
$ws.Range("A195").Value = 'Liberia'
$ws.Range("A197").Value = 'Republica del Chad'
$ws.Range("A198").Value = 'Belice'
$ws.Range("A199").Value = 'Guinea-Bisau'
$ws.Range("A203").Value = 'Papua Nueva Guinea'
$ws.Range("A204").Value = 'Timor Oriental'

# --- Updated case numbers (columns B-H) ---
$ws.Range("B4").Value = 138908
$ws.Range("C4").Value = 15330
$ws.Range("D4").Value = 4432
$ws.Range("E4").Value = 132038
$ws.Range("G4").Value = 218
$ws.Range("H4").Value = 2438
$ws.Range("B18").Value = 6280
$ws.Range("C18").Value = 625
$ws.Range("E18").Value = 5709
$ws.Range("B20").Value = 4268
$ws.Range("C20").Value = 253
$ws.Range("E20").Value = 4236
$ws.Range("B21").Value = 4256
$ws.Range("C21").Value = 352
$ws.Range("D21").Value = 6
$ws.Range("E21").Value = 4114
$ws.Range("F21").Value = 296
$ws.Range("G21").Value = 22
$ws.Range("H21").Value = 136
$ws.Range("B22").Value = 4247
$ws.Range("C22").Value = 628
$ws.Range("D22").Value = 132
$ws.Range("E22").Value = 4100
$ws.Range("F22").Value = 74
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 15
$ws.Range("B25").Value = 2775
$ws.Range("C25").Value = 144
$ws.Range("E25").Value = 2748
$ws.Range("B34").Value = 1815
$ws.Range("C34").Value = 363
$ws.Range("E34").Value = 1566
$ws.Range("B35").Value = 1597
$ws.Range("C35").Value = 102
$ws.Range("E35").Value = 1554
$ws.Range("D55").Value = 603
$ws.Range("E55").Value = 99
$ws.Range("E84").Value = 238
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 3
$ws.Range("E107").Value = 104
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 3
$ws.Range("B144").Value = 19
$ws.Range("C144").Value = 15
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 19
$ws.Range("C145").Value = 3
$ws.Range("D145").Value = 1
$ws.Range("E145").Value = 18
$ws.Range("C147").Value = 0
$ws.Range("C148").Value = 8
$ws.Range("C156").Value = 6
$ws.Range("B157").Value = 12
$ws.Range("E157").Value = 12
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 0
$ws.Range("E159").Value = 11
$ws.Range("C160").Value = 1
$ws.Range("D160").Value = 1
$ws.Range("E160").Value = 10
$ws.Range("B161").Value = 11
$ws.Range("C161").Value = 3
$ws.Range("D161").Value = 2
$ws.Range("E161").Value = 9
$ws.Range("C162").Value = 2
$ws.Range("D162").Value = 0
$ws.Range("E162").Value = 10
$ws.Range("B163").Value = 10
$ws.Range("D163").Value = 2
$ws.Range("E163").Value = 8
$ws.Range("C165").Value = 0
$ws.Range("E165").Value = 9
$ws.Range("G165").Value = 0
$ws.Range("H165").Value = 0
$ws.Range("B166").Value = 9
$ws.Range("C166").Value = 4
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 1
$ws.Range("C170").Value = 0
$ws.Range("C171").Value = 5
$ws.Range("E171").Value = 8
$ws.Range("H171").Value = 0
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 7
$ws.Range("B174").Value = 8
$ws.Range("D174").Value = 2
$ws.Range("E174").Value = 5
$ws.Range("H174").Value = 1
$ws.Range("E175").Value = 7
$ws.Range("H175").Value = 0
$ws.Range("C177").Value = 0
$ws.Range("E177").Value = 6
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 1
$ws.Range("B178").Value = 7
$ws.Range("C178").Value = 2
$ws.Range("E178").Value = 5
$ws.Range("G178").Value = 2
$ws.Range("H178").Value = 2
$ws.Range("C182").Value = 0
$ws.Range("C183").Value = 1
$ws.Range("C191").Value = 1
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 0
$ws.Range("H192").Value = 1
$ws.Range("B193").Value = 4
$ws.Range("C193").Value = 1
$ws.Range("D193").Value = 1
$ws.Range("E197").Value = 3
$ws.Range("H197").Value = 0

# --- Timestamp footer ---
$ws.Range("A1").Value = 'Datos actualizados a 29 de Marzo de 2020 a las 22:20'

